$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.496.84"
$ws.Range("E2").Value = "  -2.08%  "

$ws.Range("D3").Value = "1.749.02"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4462"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.091"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.020"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.122"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.99%  "

$ws.Range("D16").Value = "1.750.97"
$ws.Range("E16").Value = "  -1.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06381"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("E21").Value = "  -2.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.853"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("D23").Value = "27.557.99"
$ws.Range("E23").Value = "  -1.82%  "

$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.079"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("D28").Value = "1.950.31"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.088"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.080"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.653"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09017"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.535"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.70%  "

$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2084"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6339"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.943"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.205"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.378"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.739"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.717"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.145"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06859"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.54%  "
